# Weekly crime-stat refresh: roll the report forward one week and update
# the underlying 75th Precinct numbers (week-to-date / 28-day / YTD / 2-year
# columns plus all derived % change figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume number + reporting week dates -------------------------
$ws.Range("A8").Value2 = "Volume 30   Number  20"
$ws.Range("C9").Value2 = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Body: numeric value updates (counts + % change columns) --------------
    $ws.Range("M14").Value2 = -33.333333333333
    $ws.Range("N14").Value2 = -85.185185185185
    $ws.Range("F15").Value2 = 5
    $ws.Range("G15").Value2 = 5
    $ws.Range("H15").Value2 = 0
    $ws.Range("I15").Value2 = 23
    $ws.Range("J15").Value2 = 20
    $ws.Range("K15").Value2 = 15
    $ws.Range("L15").Value2 = 21.052631578947
    $ws.Range("M15").Value2 = -11.538461538461
    $ws.Range("N15").Value2 = -43.902439024390
    $ws.Range("C16").Value2 = 13
    $ws.Range("D16").Value2 = 18
    $ws.Range("E16").Value2 = -27.777777777777
    $ws.Range("F16").Value2 = 62
    $ws.Range("G16").Value2 = 71
    $ws.Range("H16").Value2 = -12.676056338028
    $ws.Range("I16").Value2 = 250
    $ws.Range("J16").Value2 = 313
    $ws.Range("K16").Value2 = -20.127795527156
    $ws.Range("L16").Value2 = 38.888888888888
    $ws.Range("M16").Value2 = -10.714285714285
    $ws.Range("N16").Value2 = -77.189781021897
    $ws.Range("C17").Value2 = 22
    $ws.Range("D17").Value2 = 22
    $ws.Range("E17").Value2 = 0
    $ws.Range("F17").Value2 = 77
    $ws.Range("G17").Value2 = 88
    $ws.Range("H17").Value2 = -12.5
    $ws.Range("I17").Value2 = 391
    $ws.Range("J17").Value2 = 425
    $ws.Range("K17").Value2 = -8
    $ws.Range("L17").Value2 = 37.192982456140
    $ws.Range("M17").Value2 = 42.181818181818
    $ws.Range("N17").Value2 = -28.388278388278
    $ws.Range("C18").Value2 = 6
    $ws.Range("D18").Value2 = 6
    $ws.Range("E18").Value2 = 0
    $ws.Range("F18").Value2 = 21
    $ws.Range("H18").Value2 = -40
    $ws.Range("I18").Value2 = 156
    $ws.Range("J18").Value2 = 164
    $ws.Range("K18").Value2 = -4.878048780487
    $ws.Range("L18").Value2 = 10.638297872340
    $ws.Range("M18").Value2 = -8.771929824561
    $ws.Range("N18").Value2 = -76.991150442477
    $ws.Range("C19").Value2 = 20
    $ws.Range("D19").Value2 = 24
    $ws.Range("E19").Value2 = -16.666666666666
    $ws.Range("F19").Value2 = 72
    $ws.Range("G19").Value2 = 108
    $ws.Range("H19").Value2 = -33.333333333333
    $ws.Range("I19").Value2 = 401
    $ws.Range("J19").Value2 = 501
    $ws.Range("K19").Value2 = -19.960079840319
    $ws.Range("L19").Value2 = 3.084832904884
    $ws.Range("M19").Value2 = 62.348178137651
    $ws.Range("N19").Value2 = 23.384615384615
    $ws.Range("C20").Value2 = 12
    $ws.Range("D20").Value2 = 11
    $ws.Range("E20").Value2 = 9.090909090909
    $ws.Range("F20").Value2 = 42
    $ws.Range("G20").Value2 = 53
    $ws.Range("H20").Value2 = -20.754716981132
    $ws.Range("I20").Value2 = 217
    $ws.Range("J20").Value2 = 219
    $ws.Range("K20").Value2 = -0.913242009132
    $ws.Range("L20").Value2 = 17.934782608695
    $ws.Range("M20").Value2 = 126.041666666667
    $ws.Range("N20").Value2 = -80.164533820840
    $ws.Range("C21").Value2 = 75
    $ws.Range("D21").Value2 = 83
    $ws.Range("E21").Value2 = -9.638554216867
    $ws.Range("F21").Value2 = 282
    $ws.Range("G21").Value2 = 361
    $ws.Range("H21").Value2 = -21.883656509695
    $ws.Range("I21").Value2 = 1446
    $ws.Range("J21").Value2 = 1648
    $ws.Range("K21").Value2 = -12.257281553398
    $ws.Range("L21").Value2 = 19.900497512437
    $ws.Range("M21").Value2 = 30.623306233062
    $ws.Range("N21").Value2 = -62.284820031298
    $ws.Range("F22").Value2 = 2
    $ws.Range("G22").Value2 = 6
    $ws.Range("H22").Value2 = -66.666666666666
    $ws.Range("J22").Value2 = 43
    $ws.Range("K22").Value2 = -58.139534883720
    $ws.Range("M22").Value2 = -35.714285714285
    $ws.Range("C23").Value2 = 4
    $ws.Range("E23").Value2 = -42.857142857142
    $ws.Range("G23").Value2 = 32
    $ws.Range("H23").Value2 = -50
    $ws.Range("I23").Value2 = 136
    $ws.Range("J23").Value2 = 136
    $ws.Range("K23").Value2 = 0
    $ws.Range("L23").Value2 = 32.038834951456
    $ws.Range("M23").Value2 = 102.985074626866
    $ws.Range("C24").Value2 = 45
    $ws.Range("D24").Value2 = 70
    $ws.Range("E24").Value2 = -35.714285714285
    $ws.Range("F24").Value2 = 198
    $ws.Range("G24").Value2 = 255
    $ws.Range("H24").Value2 = -22.352941176470
    $ws.Range("I24").Value2 = 888
    $ws.Range("J24").Value2 = 1115
    $ws.Range("K24").Value2 = -20.358744394618
    $ws.Range("L24").Value2 = 7.376058041112
    $ws.Range("M24").Value2 = 61.161524500907
    $ws.Range("C25").Value2 = 27
    $ws.Range("D25").Value2 = 28
    $ws.Range("E25").Value2 = -3.571428571428
    $ws.Range("F25").Value2 = 116
    $ws.Range("G25").Value2 = 113
    $ws.Range("H25").Value2 = 2.654867256637
    $ws.Range("I25").Value2 = 455
    $ws.Range("J25").Value2 = 494
    $ws.Range("K25").Value2 = -7.894736842105
    $ws.Range("L25").Value2 = 33.823529411764
    $ws.Range("M25").Value2 = -32.592592592592
    $ws.Range("C26").Value2 = 3
    $ws.Range("D26").Value2 = 4
    $ws.Range("E26").Value2 = -25
    $ws.Range("F26").Value2 = 8
    $ws.Range("G26").Value2 = 7
    $ws.Range("H26").Value2 = 14.285714285714
    $ws.Range("I26").Value2 = 31
    $ws.Range("J26").Value2 = 31
    $ws.Range("K26").Value2 = 0
    $ws.Range("L26").Value2 = -11.428571428571
    $ws.Range("C27").Value2 = 5
    $ws.Range("D27").Value2 = 2
    $ws.Range("E27").Value2 = 150
    $ws.Range("F27").Value2 = 15
    $ws.Range("G27").Value2 = 10
    $ws.Range("H27").Value2 = 50
    $ws.Range("I27").Value2 = 46
    $ws.Range("J27").Value2 = 46
    $ws.Range("K27").Value2 = 0
    $ws.Range("L27").Value2 = 9.523809523809
    $ws.Range("D28").Value2 = 3
    $ws.Range("E28").Value2 = -100
    $ws.Range("G28").Value2 = 8
    $ws.Range("H28").Value2 = -12.5
    $ws.Range("J28").Value2 = 31
    $ws.Range("K28").Value2 = 3.225806451612
    $ws.Range("L28").Value2 = 23.076923076923
    $ws.Range("M28").Value2 = -5.882352941176
    $ws.Range("N28").Value2 = -78.231292517006
    $ws.Range("D29").Value2 = 3
    $ws.Range("E29").Value2 = -100
    $ws.Range("G29").Value2 = 7
    $ws.Range("H29").Value2 = -14.285714285714
    $ws.Range("J29").Value2 = 28
    $ws.Range("K29").Value2 = -3.571428571428
    $ws.Range("L29").Value2 = 8
    $ws.Range("M29").Value2 = 0
    $ws.Range("N29").Value2 = -80.147058823529

# --- Cells that flip from a number to the "0" / "***.*" text placeholders -
# (these reuse the workbook's existing text styling, so copy format from a
# neighboring cell that already carries that look instead of leaving the
# quote-prefixed "typed as text" style COM applies by default)

$ws.Range("C28").Value2 = "'0"
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

$ws.Range("C29").Value2 = "'0"
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null

$ws.Range("D30").Value2 = "'0"
$ws.Range("C22").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null

$ws.Range("E30").Value2 = "'***.*"
$ws.Range("M26").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
